$wb = $excel.ActiveWorkbook

# ---- ALC ----
$ws = $wb.Worksheets.Item("ALC")
# Row 21
$ws.Range("H21").Value = 45009.5
$ws.Range("I21").Value = 45009.5
$ws.Range("K21").Value = 45009.5
$ws.Range("M21").Value = -44541.5
# Row 23
$ws.Range("H23").Value = 45009.5
$ws.Range("I23").Value = 45009.5
$ws.Range("K23").Value = 45009.5
$ws.Range("M23").Value = -44775.5
# Row 107
$ws.Range("H107").Value = 1580.3043
$ws.Range("I107").Value = 1002.9375
$ws.Range("J107").Value = 2900
$ws.Range("K107").Value = 1002.9375
$ws.Range("L107").Value = 2900
$ws.Range("M107").Value = 917.0625
$ws.Range("N107").Value = -6740
# Row 112
$ws.Range("H112").Value = 4238452.5
$ws.Range("I112").Value = 2778
$ws.Range("J112").Value = 4630645
$ws.Range("K112").Value = 8334
$ws.Range("L112").Value = 13891935
$ws.Range("M112").Value = -7226
$ws.Range("N112").Value = -13894151
# Row 137
$ws.Range("H137").Value = 2043019.2
$ws.Range("I137").Value = 3450723.8
$ws.Range("J137").Value = 1847.6
$ws.Range("K137").Value = 10352171.4
$ws.Range("L137").Value = 5542.799999999999
$ws.Range("M137").Value = -10349621.4
$ws.Range("N137").Value = -10642.8
# Row 138
$ws.Range("H138").Value = 4658.8374
$ws.Range("I138").Value = 2793.2222
$ws.Range("J138").Value = 6002.08
$ws.Range("K138").Value = 8379.6666
$ws.Range("L138").Value = 18006.24
$ws.Range("M138").Value = -3239.6666
$ws.Range("N138").Value = -28286.24
# Row 141
$ws.Range("H141").Value = 538516.7
$ws.Range("I141").Value = 1447.3529
$ws.Range("K141").Value = 4342.0587
$ws.Range("M141").Value = 837.9412999999995

# ---- ARM ----
$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 4224.9404
$ws.Range("I32").Value = 2897.1428
$ws.Range("J32").Value = 10863.929
$ws.Range("K32").Value = 2897.1428
$ws.Range("L32").Value = 10863.929
$ws.Range("M32").Value = -2610.1428
$ws.Range("N32").Value = -11437.929
# Row 97
$ws.Range("H97").Value = 650.8421
$ws.Range("I97").Value = 636.3889
$ws.Range("J97").Value = 911
$ws.Range("K97").Value = 636.3889
$ws.Range("L97").Value = 911
$ws.Range("M97").Value = -140.3889
$ws.Range("N97").Value = -1903
# Row 110
$ws.Range("H110").Value = 1656
$ws.Range("I110").Value = 682.9286
$ws.Range("J110").Value = 3602.1428
$ws.Range("K110").Value = 682.9286
$ws.Range("L110").Value = 3602.1428
$ws.Range("M110").Value = 1362.0714
$ws.Range("N110").Value = -7692.1428
# Row 139
$ws.Range("H139").Value = 24800
$ws.Range("J139").Value = 24800
$ws.Range("L139").Value = 24800
$ws.Range("N139").Value = -35080

# ---- BSM ----
$ws = $wb.Worksheets.Item("BSM")
# Row 134
$ws.Range("H134").Value = 4778.6924
$ws.Range("I134").Value = 4849.7646
$ws.Range("J134").Value = 4644.4443
$ws.Range("K134").Value = 14549.2938
$ws.Range("L134").Value = 13933.3329
$ws.Range("M134").Value = -12014.2938
$ws.Range("N134").Value = -19003.3329

# ---- CRP ----
$ws = $wb.Worksheets.Item("CRP")
# Row 16
$ws.Range("H16").Value = 2931
$ws.Range("I16").Value = 2251.6667
$ws.Range("J16").Value = 3950
$ws.Range("K16").Value = 2251.6667
$ws.Range("L16").Value = 3950
$ws.Range("M16").Value = -1964.6667
$ws.Range("N16").Value = -4524
# Row 31
$ws.Range("H31").Value = 1371981.1
$ws.Range("I31").Value = 2382222.8
$ws.Range("J31").Value = 3266.4194
$ws.Range("K31").Value = 2382222.8
$ws.Range("L31").Value = 3266.4194
$ws.Range("M31").Value = -2381927.8
$ws.Range("N31").Value = -3856.4194
# Row 34
$ws.Range("H34").Value = 1371981.1
$ws.Range("I34").Value = 2382222.8
$ws.Range("J34").Value = 3266.4194
$ws.Range("K34").Value = 2382222.8
$ws.Range("L34").Value = 3266.4194
$ws.Range("M34").Value = -2382020.8
$ws.Range("N34").Value = -3670.4194
# Row 99
$ws.Range("H99").Value = 3239.8
$ws.Range("I99").Value = 1300
$ws.Range("J99").Value = 4533
$ws.Range("K99").Value = 1300
$ws.Range("L99").Value = 4533
$ws.Range("M99").Value = 198
$ws.Range("N99").Value = -7529
# Row 107
$ws.Range("H107").Value = 1960.8695
$ws.Range("I107").Value = 473.72726
$ws.Range("J107").Value = 3324.0833
$ws.Range("K107").Value = 473.72726
$ws.Range("L107").Value = 3324.0833
$ws.Range("M107").Value = 1446.27274
$ws.Range("N107").Value = -7164.0833
# Row 113
$ws.Range("H113").Value = 2931
$ws.Range("I113").Value = 2251.6667
$ws.Range("J113").Value = 3950
$ws.Range("K113").Value = 2251.6667
$ws.Range("L113").Value = 3950
$ws.Range("M113").Value = -81.66670000000022
$ws.Range("N113").Value = -8290
# Row 126
$ws.Range("H126").Value = 3239.8
$ws.Range("I126").Value = 1300
$ws.Range("J126").Value = 4533
$ws.Range("K126").Value = 3900
$ws.Range("L126").Value = 13599
$ws.Range("M126").Value = -1430
$ws.Range("N126").Value = -18539
# Row 132
$ws.Range("H132").Value = 3132.8206
$ws.Range("I132").Value = 2012.2609
$ws.Range("K132").Value = 6036.7827
$ws.Range("M132").Value = -3506.7827
# Row 134
$ws.Range("H134").Value = 1569.7358
$ws.Range("I134").Value = 1289.7142
$ws.Range("J134").Value = 5000
$ws.Range("K134").Value = 3869.1426
$ws.Range("L134").Value = 15000
$ws.Range("M134").Value = -1334.1426
$ws.Range("N134").Value = -20070

# ---- CUL ----
$ws = $wb.Worksheets.Item("CUL")
# Row 68
$ws.Range("H68").Value = 1849.6571
$ws.Range("J68").Value = 3075.9143
$ws.Range("L68").Value = 9227.742899999999
$ws.Range("N68").Value = -10849.7429
# Row 71
$ws.Range("H71").Value = 1849.6571
$ws.Range("J71").Value = 3075.9143
$ws.Range("L71").Value = 27683.2287
$ws.Range("N71").Value = -35795.2287
# Row 92
$ws.Range("H92").Value = 2755.7144
$ws.Range("I92").Value = 1002
$ws.Range("J92").Value = 3048
$ws.Range("K92").Value = 3006
$ws.Range("L92").Value = 9144
$ws.Range("M92").Value = -1758
$ws.Range("N92").Value = -11640
# Row 107
$ws.Range("H107").Value = 849.9661
$ws.Range("J107").Value = 1139.8
$ws.Range("L107").Value = 3419.4
$ws.Range("N107").Value = -7259.4
# Row 113
$ws.Range("H113").Value = 619.82855
$ws.Range("J113").Value = 900.9286
$ws.Range("L113").Value = 2702.7858
$ws.Range("N113").Value = -7042.7858
# Row 131
$ws.Range("H131").Value = 1649.5143
$ws.Range("I131").Value = 2663.3333
$ws.Range("J131").Value = 1120.5652
$ws.Range("K131").Value = 7989.999899999999
$ws.Range("L131").Value = 3361.6956
$ws.Range("M131").Value = -2949.999899999999
$ws.Range("N131").Value = -13441.6956
# Row 132
$ws.Range("H132").Value = 2365.9524
$ws.Range("J132").Value = 2274.25
$ws.Range("L132").Value = 20468.25
$ws.Range("N132").Value = -25528.25

# ---- GSM ----
$ws = $wb.Worksheets.Item("GSM")
# Row 132
$ws.Range("H132").Value = 3603.5676
$ws.Range("I132").Value = 2573.56
$ws.Range("K132").Value = 7720.68
$ws.Range("M132").Value = -5190.68

# ---- LTW ----
$ws = $wb.Worksheets.Item("LTW")
# Row 22
$ws.Range("H22").Value = 250002880
$ws.Range("I22").Value = 333333820
$ws.Range("J22").Value = 10000
$ws.Range("K22").Value = 333333820
$ws.Range("L22").Value = 10000
$ws.Range("M22").Value = -333333525
$ws.Range("N22").Value = -10590
# Row 27
$ws.Range("H27").Value = 250002880
$ws.Range("I27").Value = 333333820
$ws.Range("J27").Value = 10000
$ws.Range("K27").Value = 333333820
$ws.Range("L27").Value = 10000
$ws.Range("M27").Value = -333333713
$ws.Range("N27").Value = -10214
# Row 61
$ws.Range("H61").Value = 100003736
$ws.Range("I61").Value = 125002170
$ws.Range("J61").Value = 10000
$ws.Range("K61").Value = 125002170
$ws.Range("L61").Value = 10000
$ws.Range("M61").Value = -125001968
$ws.Range("N61").Value = -10404
# Row 113
$ws.Range("H113").Value = 100003736
$ws.Range("I113").Value = 125002170
$ws.Range("J113").Value = 10000
$ws.Range("K113").Value = 125002170
$ws.Range("L113").Value = 10000
$ws.Range("M113").Value = -125000000
$ws.Range("N113").Value = -14340
# Row 138
$ws.Range("H138").Value = 50000
$ws.Range("J138").Value = 50000
$ws.Range("L138").Value = 50000
$ws.Range("N138").Value = -60280

# ---- WVR ----
$ws = $wb.Worksheets.Item("WVR")
# Row 107
$ws.Range("H107").Value = 2387.5715
$ws.Range("I107").Value = 1000
$ws.Range("J107").Value = 2494.3076
$ws.Range("K107").Value = 3000
$ws.Range("L107").Value = 7482.9228
$ws.Range("M107").Value = -1080
$ws.Range("N107").Value = -11322.9228
# Row 113
$ws.Range("H113").Value = 2555.1428
$ws.Range("I113").Value = 633.3333
$ws.Range("J113").Value = 3996.5
$ws.Range("K113").Value = 1899.9999
$ws.Range("L113").Value = 11989.5
$ws.Range("M113").Value = 270.0001
$ws.Range("N113").Value = -16329.5
